# Generate Report for Handback
#
# This mirrors a "handback" report run: the localized xliff files came back
# in sync with en-US, so:
#   - the human-readable status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown,
#   - the per-language detail sheets (zh-cn / de-de) get their "Latest
#     Target File" / "Latest Handback File" / "Latest Handback DateTime"
#     columns (I/J/K) filled in for the two tracked source files,
#   - the "Latest Target File" cell becomes a hyperlink back to the source
#     file on GitHub, matching the existing hyperlink already shown in
#     column A.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa2df44b16570f24ec2a539b01a93673f64ac5c3/e2e/ea75a6df-8348-434a-aa74-dd50eb0270f8.md"
$targetDisplay = "ea75a6df-8348-434a-aa74-dd50eb0270f8.md"

# ---------------------------------------------------------------------
# Overview sheet: status column for both tracked files (E/F, rows 2-3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Per-language detail sheets
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; XlfTarget = "ea75a6df-8348-434a-aa74-dd50eb0270f8.5d8bbd2f4b3f591a893b0c230c9634346d93c616.zh-cn.xlf"; HandbackDate = "2016-08-19 04:57:08" },
    @{ Sheet = "de-de"; XlfTarget = "ea75a6df-8348-434a-aa74-dd50eb0270f8.5d8bbd2f4b3f591a893b0c230c9634346d93c616.de-de.xlf"; HandbackDate = "2016-08-19 04:57:15" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for both rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Handback File (J) and Latest Handback DateTime (K)
    $ws.Range("J2").Value = $lang.XlfTarget
    $ws.Range("K2").Value = $lang.HandbackDate
    $ws.Range("J3").Value = $lang.XlfTarget
    $ws.Range("K3").Value = $lang.HandbackDate

    # Latest Target File (I) becomes a hyperlink to the source file,
    # matching the hyperlink already used in column A.
    $ws.Range("I2").Value = $targetDisplay
    $ws.Range("I3").Value = $targetDisplay
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, "", "", $targetDisplay)
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, "", "", $targetDisplay)

    # Column widths: Status (C) grows with the longer text, Latest
    # Target File / Latest Handback File (I/J) widen to fit filenames.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
